# Auto-generated edit script applying the Pandaemonium_Profits.xlsx diff
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) for several leves
# across the ALC, ARM, BSM, CRP, CUL, GSM and LTW sheets.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 725.5714
$ws.Range("I92").Value = 732.0909
$ws.Range("J92").Value = 701.6667
$ws.Range("K92").Value = 732.0909
$ws.Range("L92").Value = 701.6667
$ws.Range("M92").Value = 515.9091
$ws.Range("N92").Value = -3197.6667
$ws.Range("H112").Value = 5799.9375
$ws.Range("J112").Value = 1696.0769
$ws.Range("L112").Value = 5088.2307
$ws.Range("N112").Value = -7304.2307
$ws.Range("H129").Value = 928.8043
$ws.Range("I129").Value = 348.5
$ws.Range("J129").Value = 955.1818
$ws.Range("K129").Value = 1045.5
$ws.Range("L129").Value = 2865.5454
$ws.Range("M129").Value = 3954.5
$ws.Range("N129").Value = -12865.5454

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("H5").Value = 263
$ws.Range("I5").Value = 258.8889
$ws.Range("K5").Value = 258.8889
$ws.Range("M5").Value = -146.8889
$ws.Range("H63").Value = 2908.9375
$ws.Range("I63").Value = 2519.8333
$ws.Range("J63").Value = 4076.25
$ws.Range("K63").Value = 2519.8333
$ws.Range("L63").Value = 4076.25
$ws.Range("M63").Value = -1833.8333
$ws.Range("N63").Value = -5448.25
$ws.Range("H66").Value = 2908.9375
$ws.Range("I66").Value = 2519.8333
$ws.Range("J66").Value = 4076.25
$ws.Range("K66").Value = 12599.1665
$ws.Range("L66").Value = 20381.25
$ws.Range("M66").Value = -9167.166499999999
$ws.Range("N66").Value = -27245.25
$ws.Range("H88").Value = 6777.5835
$ws.Range("I88").Value = 18133.334
$ws.Range("J88").Value = 2992.3333
$ws.Range("K88").Value = 18133.334
$ws.Range("L88").Value = 2992.3333
$ws.Range("M88").Value = -17727.334
$ws.Range("N88").Value = -3804.3333
$ws.Range("H91").Value = 6777.5835
$ws.Range("I91").Value = 18133.334
$ws.Range("J91").Value = 2992.3333
$ws.Range("K91").Value = 18133.334
$ws.Range("L91").Value = 2992.3333
$ws.Range("M91").Value = -16729.334
$ws.Range("N91").Value = -5800.3333
$ws.Range("H132").Value = 2580.4666
$ws.Range("I132").Value = 2533.1482
$ws.Range("J132").Value = 2651.4443
$ws.Range("K132").Value = 7599.444600000001
$ws.Range("L132").Value = 7954.3329
$ws.Range("M132").Value = -5069.444600000001
$ws.Range("N132").Value = -13014.3329
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 263
$ws.Range("I4").Value = 258.8889
$ws.Range("K4").Value = 258.8889
$ws.Range("M4").Value = -143.8889
$ws.Range("H82").Value = 16867.334
$ws.Range("I82").Value = 5159.5
$ws.Range("J82").Value = 40283
$ws.Range("K82").Value = 5159.5
$ws.Range("L82").Value = 40283
$ws.Range("M82").Value = -4776.5
$ws.Range("N82").Value = -41049
$ws.Range("H85").Value = 16867.334
$ws.Range("I85").Value = 5159.5
$ws.Range("J85").Value = 40283
$ws.Range("K85").Value = 5159.5
$ws.Range("L85").Value = 40283
$ws.Range("M85").Value = -3833.5
$ws.Range("N85").Value = -42935
$ws.Range("H86").Value = 2026.0238
$ws.Range("I86").Value = 1877.6875
$ws.Range("J86").Value = 2500.7
$ws.Range("K86").Value = 1877.6875
$ws.Range("L86").Value = 2500.7
$ws.Range("M86").Value = -754.6875
$ws.Range("N86").Value = -4746.7
$ws.Range("H89").Value = 2026.0238
$ws.Range("I89").Value = 1877.6875
$ws.Range("J89").Value = 2500.7
$ws.Range("K89").Value = 9388.4375
$ws.Range("L89").Value = 12503.5
$ws.Range("M89").Value = -3772.4375
$ws.Range("N89").Value = -23735.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2181.7273
$ws.Range("I99").Value = 1699.8334
$ws.Range("K99").Value = 1699.8334
$ws.Range("M99").Value = -201.8334
$ws.Range("H107").Value = 1090.3636
$ws.Range("I107").Value = 1200.091
$ws.Range("J107").Value = 980.63635
$ws.Range("K107").Value = 1200.091
$ws.Range("L107").Value = 980.63635
$ws.Range("M107").Value = 719.9090000000001
$ws.Range("N107").Value = -4820.63635
$ws.Range("H126").Value = 2181.7273
$ws.Range("I126").Value = 1699.8334
$ws.Range("K126").Value = 5099.5002
$ws.Range("M126").Value = -2629.5002
$ws.Range("H132").Value = 2908.5745
$ws.Range("I132").Value = 2661
$ws.Range("J132").Value = 3388.25
$ws.Range("K132").Value = 7983
$ws.Range("L132").Value = 10164.75
$ws.Range("M132").Value = -5453
$ws.Range("N132").Value = -15224.75

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 156.2
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 170.25
$ws.Range("K7").Value = 300
$ws.Range("L7").Value = 510.75
$ws.Range("M7").Value = -188
$ws.Range("N7").Value = -734.75
$ws.Range("H9").Value = 36253.617
$ws.Range("J9").Value = 36253.617
$ws.Range("L9").Value = 108760.851
$ws.Range("N9").Value = -109208.851
$ws.Range("H17").Value = 5002
$ws.Range("J17").Value = 5002
$ws.Range("L17").Value = 15006
$ws.Range("N17").Value = -15344
$ws.Range("H55").Value = 8828.571
$ws.Range("J55").Value = 9983.333000000001
$ws.Range("L55").Value = 29949.999
$ws.Range("N55").Value = -30303.999
$ws.Range("H80").Value = 2924.75
$ws.Range("I80").Value = 2899
$ws.Range("J80").Value = 2933.3333
$ws.Range("K80").Value = 8697
$ws.Range("L80").Value = 8799.999899999999
$ws.Range("M80").Value = -7761
$ws.Range("N80").Value = -10671.9999
$ws.Range("H83").Value = 2924.75
$ws.Range("I83").Value = 2899
$ws.Range("J83").Value = 2933.3333
$ws.Range("K83").Value = 26091
$ws.Range("L83").Value = 26399.9997
$ws.Range("M83").Value = -21411
$ws.Range("N83").Value = -35759.9997
$ws.Range("H117").Value = 3086
$ws.Range("J117").Value = 3086
$ws.Range("L117").Value = 9258
$ws.Range("N117").Value = -16142
$ws.Range("H123").Value = 430
$ws.Range("I123").Value = 430
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 1290
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = 1160
$ws.Range("H131").Value = 1413.2174
$ws.Range("I131").Value = 2435.7144
$ws.Range("J131").Value = 1229.6923
$ws.Range("K131").Value = 7307.1432
$ws.Range("L131").Value = 3689.0769
$ws.Range("M131").Value = -2267.1432
$ws.Range("N131").Value = -13769.0769
$ws.Range("N123").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5681.4883
$ws.Range("I70").Value = 5446.3076
$ws.Range("J70").Value = 6041.1763
$ws.Range("K70").Value = 5446.3076
$ws.Range("L70").Value = 6041.1763
$ws.Range("M70").Value = -5176.3076
$ws.Range("N70").Value = -6581.1763
$ws.Range("H73").Value = 5681.4883
$ws.Range("I73").Value = 5446.3076
$ws.Range("J73").Value = 6041.1763
$ws.Range("K73").Value = 5446.3076
$ws.Range("L73").Value = 6041.1763
$ws.Range("M73").Value = -4510.3076
$ws.Range("N73").Value = -7913.1763

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 430.4
$ws.Range("I22").Value = 459.4
$ws.Range("J22").Value = 401.4
$ws.Range("K22").Value = 459.4
$ws.Range("L22").Value = 401.4
$ws.Range("M22").Value = -164.4
$ws.Range("N22").Value = -991.4
$ws.Range("H27").Value = 430.4
$ws.Range("I27").Value = 459.4
$ws.Range("J27").Value = 401.4
$ws.Range("K27").Value = 459.4
$ws.Range("L27").Value = 401.4
$ws.Range("M27").Value = -352.4
$ws.Range("N27").Value = -615.4
$ws.Range("H46").Value = 1400
$ws.Range("J46").Value = 2000
$ws.Range("L46").Value = 2000
$ws.Range("N46").Value = -2376
$ws.Range("H100").Value = 4410.7646
$ws.Range("I100").Value = 3798.8667
$ws.Range("K100").Value = 3798.8667
$ws.Range("M100").Value = -3257.8667
$ws.Range("H122").Value = 6065.125
$ws.Range("I122").Value = 6154.839
$ws.Range("J122").Value = 5756.1113
$ws.Range("K122").Value = 18464.517
$ws.Range("L122").Value = 17268.3339
$ws.Range("M122").Value = -16014.517
$ws.Range("N122").Value = -22168.3339
